$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Cash
$ws.Range("C2").Value = 0.03398908557859916
$ws.Range("D2").Value = 0.03398899316178236

# Row 3 - EU Flot
$ws.Range("C3").Value = 0.03326733264542253
$ws.Range("D3").Value = 0.03326733196496505

# Row 4 - EU Equity
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0.09196051005819827
$ws.Range("D4").Value = 0.09196052048333501

# Row 5 - US Equity
$ws.Range("C5").Value = 0.09359360785550738
$ws.Range("D5").Value = 0.09359362526776671

# Row 6 - Greek Gov
$ws.Range("C6").Value = 0.6494120724745848
$ws.Range("D6").Value = 0.64941204950374

# Row 7 - EU Corps
$ws.Range("C7").Value = 0.03665620485808728
$ws.Range("D7").Value = 0.03665631947025359

# Row 8 - EU Gov
$ws.Range("C8").Value = 0.06112118652960059
$ws.Range("D8").Value = 0.06112116014815723
